# Stuller Classic Bands Mens - parse update
# - Rename the worksheet from the generic "Sheet1" to a descriptive name.
# - Collapse the per-variant SKU values in Model (A) / Sku (B) down to the
#   shared base model number for each width group (the individual metal /
#   karat suffix is dropped since it's already captured by other columns).
# - Normalize the capitalization of the Metal column for the Platinum and
#   Palladium variants in the 4mm and 5mm groups ("platinum "/"palladium"
#   -> "Platinum "/"Palladium").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Stuller Classic Bands Mens"

# Row ranges (inclusive) for each width group and their shared base model #.
$groups = @(
    @{ Start = 2;  End = 9;  Model = "BHRL75M" },
    @{ Start = 10; End = 16; Model = "BHRL74M" },
    @{ Start = 17; End = 24; Model = "BHRL76M" },
    @{ Start = 25; End = 31; Model = "BHRL77M" },
    @{ Start = 32; End = 38; Model = "BHRL78M" }
)

foreach ($group in $groups) {
    for ($r = $group.Start; $r -le $group.End; $r++) {
        $ws.Cells.Item($r, 1).Value = $group.Model
        $ws.Cells.Item($r, 2).Value = $group.Model
    }
}

# Capitalize Metal (column J) for the Platinum / Palladium rows of the
# 5mm (rows 8-9) and 4mm (rows 15-16) groups only.
$ws.Cells.Item(8, 10).Value = "Platinum "
$ws.Cells.Item(9, 10).Value = "Palladium"
$ws.Cells.Item(15, 10).Value = "Platinum "
$ws.Cells.Item(16, 10).Value = "Palladium"

# Restore the scrolled view state (top-left cell + active selection) that
# was recorded in the saved workbook.
$win = $excel.ActiveWindow
$win.ScrollRow = 25
$win.ScrollColumn = 1
$ws.Range("H39").Select()
